# Apply "hybrid bold + color" highlighting to quantitative impact metrics
# across the resume: wraps each metric substring (percentages, dollar
# amounts, etc.) in its own run with Bold + color 2C3E50, matching the
# target diff exactly (run-splitting is left to the host engine, which
# auto-inserts xml:space="preserve" on the surrounding text runs).

$d = $word.ActiveDocument

# wdColor value for RRGGBB hex 2C3E50 (Word packs colors as 0xBBGGRR)
$HighlightColor = 5258796

function Format-Metric {
    param($SearchRange, [string]$Needle)
    $found = $SearchRange.Find.Execute($Needle, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if ($found) {
        $SearchRange.Font.Bold = 1
        $SearchRange.Font.Color = $HighlightColor
    }
    return $found
}

function Format-MetricsInParagraph {
    param($Paragraph, [string[]]$Needles)
    if ($Paragraph -eq $null) {
        return
    }
    $paraEnd = $Paragraph.Range.End
    $cursor = $Paragraph.Range.Start
    foreach ($needle in $Needles) {
        $searchRange = $d.Range($cursor, $paraEnd)
        $found = Format-Metric $searchRange $needle
        if ($found) {
            $cursor = $searchRange.End
        }
    }
}

function Get-ParagraphContaining {
    param([string]$Marker)
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$Marker*") {
            return $p
        }
    }
    return $null
}

# Paragraph: "• Discovered systematic race coding errors ... from 23% to 64%"
$p = Get-ParagraphContaining "Discovered systematic race coding errors"
Format-MetricsInParagraph $p @("23%", "64%")

# Paragraph: "• Utilized advanced sampling methods ... ±4.2% to ±2.1% ... 71% to 87% ..."
$p = Get-ParagraphContaining "Utilized advanced sampling methods"
Format-MetricsInParagraph $p @("±4.2%", "±2.1%", "71%", "87%")

# Paragraph: "• Trigonometric algorithm for boundary estimation ... 73.5% ... $4.7M ..."
$p = Get-ParagraphContaining "Trigonometric algorithm for boundary estimation"
Format-MetricsInParagraph $p @("73.5%", "`$4.7M")

# Paragraph: "• Built real-time FEC analysis systems ... valued over $2 trillion"
$p = Get-ParagraphContaining "Built real-time FEC analysis systems"
Format-MetricsInParagraph $p @("`$2")

# Paragraph: "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
$p = Get-ParagraphContaining "Algorithmic innovation: Pioneered trigonometric"
Format-MetricsInParagraph $p @("73.5%")

# Paragraph: "• $4.7M savings enabled nonprofit access"
$p = Get-ParagraphContaining "4.7M savings enabled nonprofit access"
Format-MetricsInParagraph $p @("`$4.7M")

# Paragraph: "• 178% accuracy improvement in racial classification algorithms"
$p = Get-ParagraphContaining "178% accuracy improvement in racial classification"
Format-MetricsInParagraph $p @("178%")
